# Update the "Pressure" column (B2:B49) with the new measurements and
# strip the ad-hoc "Обычный 2" cell style (centered / wrapped 12pt font)
# that used to be applied to those cells, matching the author's re-export
# of the lab plot data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(750,749,750,752,753,752,751,754,753,754,755,755,753,752,753,753,751,752,755,757,755,755,757,758,755,753,753,749,745,743,743,744,749,754,752,751,750,746,750,757,757,752,750,747,740,738,735,738)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# The cells no longer carry the custom "Обычный 2" style - clear direct
# formatting back to the default style.
$ws.Range("B2:B49").ClearFormats()

# Drop the now-unused custom cell style from the workbook style table.
for ($i = $wb.Styles.Count; $i -ge 1; $i--) {
    $style = $wb.Styles.Item($i)
    if ($style.Name -ne "Normal" -and $style.Name -ne "Обычный") {
        $style.Delete()
    }
}

# Move the viewport / selection the way the author left it.
$ws.Range("H38").Select()
